$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.8683625768187075

$ws.Range("H3").Value = 0.8693849133062217

$ws.Range("H4").Value = 0.8694219203650375

$ws.Range("A5").Value = 6
$ws.Range("B5").Value = 'XGBRegressor'
$ws.Range("C5").Value = '{''learning_rate'': 0.1, ''max_depth'': 3, ''n_estimators'': 100, ''random_state'': 42, ''subsample'': 0.8}'
$ws.Range("D5").Value = 0.869636058807373
$ws.Range("E5").Value = 4.366933994293213
$ws.Range("F5").Value = 5.632278085654206
$ws.Range("G5").Value = 31.72255643414061
$ws.Range("H5").Value = 0.8529651999473572

$ws.Range("A6").Value = 7
$ws.Range("B6").Value = 'CatBoostRegressor'
$ws.Range("C6").Value = '{''depth'': 5, ''iterations'': 1000, ''l2_leaf_reg'': 3, ''learning_rate'': 0.01}'
$ws.Range("D6").Value = 0.8659973357237966
$ws.Range("E6").Value = 4.367474193136229
$ws.Range("F6").Value = 5.710340963612286
$ws.Range("G6").Value = 32.60799392070849
$ws.Range("H6").Value = 0.8609961825435398

$ws.Range("C7").Value = '{''max_depth'': 10, ''min_samples_leaf'': 2, ''min_samples_split'': 10, ''n_estimators'': 100, ''random_state'': 42}'
$ws.Range("D7").Value = 0.850382474334519
$ws.Range("E7").Value = 4.595319420439647
$ws.Range("F7").Value = 6.03387846309462
$ws.Range("G7").Value = 36.4076893073971
$ws.Range("H7").Value = 0.8405518330742539

$ws.Range("H8").Value = 0.8224052948218432

$ws.Range("H9").Value = 0.8057104066076171

$ws.Range("H10").Value = 0.6751133216615145

